$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 8.5
$ws.Range("R2").Value = 2.2
$ws.Range("S2").Value = 1.62
$ws.Range("V2").Value = 9.5
$ws.Range("X2").Value = 15
$ws.Range("AA2").Value = 8.5
$ws.Range("AB2").Value = 23
$ws.Range("AE2").Value = 17
$ws.Range("AG3").Value = 12
$ws.Range("G8").Value = 2.35
$ws.Range("H8").Value = 2.88
$ws.Range("I8").Value = 3.5
$ws.Range("N8").Value = 2.88
$ws.Range("O8").Value = 1.4
$ws.Range("R8").Value = 2.25
$ws.Range("S8").Value = 1.57
$ws.Range("U8").Value = 9.5
$ws.Range("V8").Value = 11
$ws.Range("W8").Value = 23
$ws.Range("X8").Value = 23
$ws.Range("Z8").Value = 5.5
$ws.Range("AA8").Value = 6
$ws.Range("AE8").Value = 7.5
$ws.Range("AF8").Value = 15
$ws.Range("AG8").Value = 13
$ws.Range("AI8").Value = 34
$ws.Range("G9").Value = 4.1
$ws.Range("H9").Value = 2.88
$ws.Range("J9").Value = 1.17
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 1.67
$ws.Range("M9").Value = 2.1
$ws.Range("N9").Value = 3.4
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 1.75
$ws.Range("Q9").Value = 2.05
$ws.Range("R9").Value = 2.63
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 7.5
$ws.Range("U9").Value = 17
$ws.Range("V9").Value = 17
$ws.Range("W9").Value = 41
$ws.Range("X9").Value = 41
$ws.Range("Y9").Value = 67
$ws.Range("Z9").Value = 5
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 101
$ws.Range("AE9").Value = 5
$ws.Range("AF9").Value = 8.5
$ws.Range("AG9").Value = 11
$ws.Range("AI9").Value = 26
$ws.Range("I10").Value = 4
$ws.Range("K10").Value = 13
$ws.Range("N10").Value = 1.75
$ws.Range("O10").Value = 2.05
$ws.Range("T10").Value = 8
$ws.Range("W10").Value = 15
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 4.33
$ws.Range("J11").Value = 1.04
$ws.Range("K11").Value = 13
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4
$ws.Range("N11").Value = 1.73
$ws.Range("O11").Value = 2.08
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.7
$ws.Range("S11").Value = 2.05
$ws.Range("T11").Value = 8.5
$ws.Range("U11").Value = 9
$ws.Range("Y11").Value = 23
$ws.Range("Z11").Value = 13
$ws.Range("AB11").Value = 15
$ws.Range("AG11").Value = 15
$ws.Range("AI11").Value = 34
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 1.67
$ws.Range("J13").Value = 1.08
$ws.Range("K13").Value = 7.5
$ws.Range("L13").Value = 1.4
$ws.Range("M13").Value = 2.75
$ws.Range("N13").Value = 2.35
$ws.Range("O13").Value = 1.57
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.5
$ws.Range("R13").Value = 2.2
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 11
$ws.Range("V13").Value = 19
$ws.Range("X13").Value = 51
$ws.Range("Z13").Value = 7.5
$ws.Range("AB13").Value = 21
$ws.Range("AC13").Value = 81
$ws.Range("AE13").Value = 5.5
$ws.Range("AG13").Value = 9
$ws.Range("AI13").Value = 17
$ws.Range("L14").Value = 1.17
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 1.6
$ws.Range("O14").Value = 2.3
$ws.Range("G15").Value = 2.7
$ws.Range("I15").Value = 2.63
$ws.Range("O15").Value = 1.6
$ws.Range("U15").Value = 12
$ws.Range("Y15").Value = 41
$ws.Range("Z15").Value = 7.5
$ws.Range("AD15").Value = 401
$ws.Range("AG15").Value = 11
$ws.Range("G16").Value = 3.2
$ws.Range("I16").Value = 2.2
$ws.Range("O16").Value = 1.7
$ws.Range("U16").Value = 15
$ws.Range("Y16").Value = 41
$ws.Range("Z16").Value = 8.5
$ws.Range("AG16").Value = 9.5
$ws.Range("AH16").Value = 21
$ws.Range("N17").Value = 2.08
$ws.Range("O17").Value = 1.73
$ws.Range("G18").Value = 1.91
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 4
$ws.Range("O18").Value = 1.57
$ws.Range("W18").Value = 17
$ws.Range("X18").Value = 19
$ws.Range("AE18").Value = 9
$ws.Range("AF18").Value = 19
$ws.Range("AH18").Value = 41
$ws.Range("L19").Value = 1.4
$ws.Range("M19").Value = 2.75
$ws.Range("N19").Value = 2.3
$ws.Range("O19").Value = 1.6
$ws.Range("P19").Value = 1.5
$ws.Range("Q19").Value = 2.5
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1.75
$ws.Range("AA19").Value = 6.5
$ws.Range("AB19").Value = 17
$ws.Range("AD19").Value = 451
$ws.Range("AE19").Value = 6.5
$ws.Range("AI19").Value = 19
$ws.Range("J20").Value = 1.05
$ws.Range("K20").Value = 11
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("R20").Value = 2.5
$ws.Range("S20").Value = 1.5
$ws.Range("T20").Value = 5.5
$ws.Range("X20").Value = 13
$ws.Range("Z20").Value = 9
$ws.Range("AA20").Value = 9.5
$ws.Range("N22").Value = 2
$ws.Range("O22").Value = 1.85
$ws.Range("N23").Value = 1.83
$ws.Range("O23").Value = 2.03
$ws.Range("G24").Value = 3.8
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 2.1
$ws.Range("X24").Value = 34
$ws.Range("AA24").Value = 6
$ws.Range("H26").Value = 3.5
$ws.Range("I26").Value = 1.95
$ws.Range("N26").Value = 1.93
$ws.Range("O26").Value = 1.93
$ws.Range("W26").Value = 41
$ws.Range("Z26").Value = 11
$ws.Range("AE26").Value = 8
$ws.Range("G29").Value = 1.7
$ws.Range("H29").Value = 3.9
$ws.Range("I29").Value = 4.3
$ws.Range("N29").Value = 1.62
$ws.Range("P29").Value = 1.3
$ws.Range("Q29").Value = 3.2
$ws.Range("R29").Value = 1.62
$ws.Range("S29").Value = 2.15
$ws.Range("T29").Value = 8.75
$ws.Range("U29").Value = 9.25
$ws.Range("W29").Value = 14
$ws.Range("Y29").Value = 21
$ws.Range("AA29").Value = 7.7
$ws.Range("AB29").Value = 13.5
$ws.Range("AE29").Value = 15
$ws.Range("AF29").Value = 27
$ws.Range("AG29").Value = 14
$ws.Range("AH29").Value = 70
$ws.Range("AI29").Value = 37
